# Remove the trailing sentence
#   "Two ways to convey "nothing there" won't seem so strange in Chapter 3."
# from the paragraph about JavaScript value types.
#
# Search every slide/shape/paragraph for the old text so the edit is robust
# to shape/slide index ordering.

$oldSentenceMarker = "Two ways to convey"
$newText = "In JavaScript, data is represented with values. There are four value types to convey data with: string, number, boolean, and object. Additionally, there are two value types to convey no data with: undefined and null."

$p = $ppt.ActivePresentation
$found = $false

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTextFrame -ne $false) {
            $tf = $shape.TextFrame
            if ($tf.HasText -ne $false) {
                $tr = $tf.TextRange
                $paraCount = $tr.Paragraphs().Count
                for ($k = 1; $k -le $paraCount; $k++) {
                    $para = $tr.Paragraphs($k)
                    if ($para.Text -like "*$oldSentenceMarker*") {
                        $para.Text = $newText
                        $found = $true
                    }
                }
            }
        }
    }
}

Write-Output "Replaced: $found"
